$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Copy the "Durchschnitt / nur für / Service" label block (F27:H27) down
#    to rows 29 and 30 BEFORE row 27 itself gets re-purposed for the new
#    "Standardabweichung" row (its formatting is about to be overwritten).
# ---------------------------------------------------------------------------
$ws.Range("F27:H27").Copy()
$ws.Range("F29:H29").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("F27:H27").Copy()
$ws.Range("F30:H30").PasteSpecial(-4122)   # xlPasteFormats

# ---------------------------------------------------------------------------
# 2) Turn the previously-empty rows 26 and 27 into the "Median" and
#    "Standardabweichung" summary rows by copying the formatting already
#    used for row 25 ("Durchschnitt").
# ---------------------------------------------------------------------------
$ws.Range("C25:P25").Copy()
$ws.Range("C26:P26").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("C25:P25").Copy()
$ws.Range("C27:P27").PasteSpecial(-4122)   # xlPasteFormats

# ---------------------------------------------------------------------------
# 3) Fill in the new text labels, in the same order they appear in the
#    target shared-strings table (Median, Std.-Abw., Standardabweichung).
# ---------------------------------------------------------------------------
$ws.Range("C26").Value = "Median"
$ws.Range("F30").Value = "Std.-Abw."
$ws.Range("C27").Value = "Standardabweichung"

$ws.Range("F29").Value = "Durchschnitt"
$ws.Range("G29").Value = "nur für"
$ws.Range("H29").Value = "Service"

$ws.Range("G30").Value = "nur für"
$ws.Range("H30").Value = "Service"

# ---------------------------------------------------------------------------
# 4) Formulas for the Median row (26).
# ---------------------------------------------------------------------------
$ws.Range("K26").Formula = "=MEDIAN(K7:K21)"
$ws.Range("L26").Formula = "=MEDIAN(L7:L21)"

# ---------------------------------------------------------------------------
# 5) Formulas for the Standardabweichung row (27).
# ---------------------------------------------------------------------------
$ws.Range("F27").Formula = "=ROUND(STDEV.P(F8,F10,F12,F14,F16,F18,F20,F22),1)"
$ws.Range("G27").Formula = "=ROUND(STDEV.P(G8,G10,G12,G14,G16,G18,G20,G22),1)"
$ws.Range("H27").Formula = "=ROUND(STDEV.P(H8,H10,H12,H14,H16,H18,H20,H22),1)"
$ws.Range("K27").Formula = "=ROUND(STDEV.P(K7:K21),1)"
$ws.Range("L27").Formula = "=ROUND(STDEV.P(L7:L21),1)"

# ---------------------------------------------------------------------------
# 6) Widen column C to fit the new, longer row labels.
# ---------------------------------------------------------------------------
$ws.Columns(3).ColumnWidth = 23.83

# ---------------------------------------------------------------------------
# 7) Scroll/selection: the sheet view now starts at A3 with L29 selected.
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 3
$win.ScrollColumn = 1
$ws.Range("L29").Select()
